# feat: add 2022-Q1 data
#
# 1) Duplicate the "2021-Q4" sheet (same column layout) to create the new
#    "2022-Q1" sheet, placed immediately before "总计".
# 2) Fill the new sheet with the 2022-Q1 fund holdings data (8 rows).
# 3) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q1 (count=8, value=3.09) and shift the previous rows down,
#    renumbering the index column.

$wb = $excel.ActiveWorkbook

$sourceSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet  = $wb.Worksheets.Item("总计")

# Duplicate the source sheet (keeps header/row styles identical) and place
# it right before the "总计" sheet, then rename it.
$sourceSheet.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The source sheet only has 5 data rows (rows 2-6); we need 8 (rows 2-9).
# Extend the index-column formatting (style) down to the new rows by
# copying the format of the last existing row.
$newSheet.Range("A6").Copy()
$newSheet.Range("A7:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Holdings data for 2022-Q1: index, code, name, fund size, total equity
# position, position ratio, held value (100M yuan), position rank.
$fundRows = @(
  @(0, "011230", "创金合信数字经济主题股票C", "17.18", "92.17", "6.04", "1.0377", 2),
  @(1, "166025", "中欧远见两年定期开放混合A", "48.80", "48.74", "2.06", "1.0053", 9),
  @(2, "011229", "创金合信数字经济主题股票A", "12.18", "92.17", "6.04", "0.7357", 2),
  @(3, "206002", "鹏华精选成长混合", "4.48", "92.68", "4.61", "0.2065", 8),
  @(4, "007101", "中欧远见两年定期开放混合C", "2.87", "48.74", "2.06", "0.0591", 9),
  @(5, "000066", "诺安鸿鑫混合", "0.74", "81.34", "5.25", "0.0388", 2),
  @(6, "005104", "富荣福康混合A", "0.08", "87.88", "3.00", "0.0024", 9),
  @(7, "005105", "富荣福康混合C", "0.04", "87.88", "3.00", "0.0012", 9)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
  $r = $i + 2
  $row = $fundRows[$i]
  # Column A: numeric row index
  $newSheet.Cells.Item($r, 1).Value = $row[0]
  # Columns B, D, E, F, G are stored as text (leading zeros / trailing
  # zeros must be preserved), so prefix with an apostrophe to force text.
  $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
  # Column C is a plain text name, no ambiguity with numbers.
  $newSheet.Cells.Item($r, 3).Value = $row[2]
  $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
  $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
  $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
  $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
  # Column H: numeric rank.
  $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# --- Update the "总计" summary sheet -------------------------------------

$ws = $wb.Worksheets.Item("总计")

# Remember the existing quarterly rows (currently rows 2-5) before
# overwriting them.
$existing = @()
for ($r = 2; $r -le 5; $r++) {
  $existing += ,@($ws.Cells.Item($r, 2).Value(), $ws.Cells.Item($r, 3).Value(), $ws.Cells.Item($r, 4).Value())
}

# New first row: the 2022-Q1 totals.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 8
$ws.Cells.Item(2, 4).Value = 3.09

# Re-write the previously existing rows, shifted down by one, with the
# index column renumbered accordingly.
for ($i = 0; $i -lt $existing.Length; $i++) {
  $r = $i + 3
  $ws.Cells.Item($r, 1).Value = $i + 1
  $ws.Cells.Item($r, 2).Value = $existing[$i][0]
  $ws.Cells.Item($r, 3).Value = $existing[$i][1]
  $ws.Cells.Item($r, 4).Value = $existing[$i][2]
}

# Row 6 is brand new (sheet previously only spanned to row 5); give its
# index cell (A6) the same style as the row above it (A5).
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0
